# Actualización desde MV -datos-
# Adds a new publication-vintage column (BH, "Agosto.2021") that mostly repeats
# the previous vintage (BG) except for the most recent observation, and adds a
# new observation row (75, period "01-04-2021") for that new vintage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header: BH1 = "Agosto.2021" (copy the bold/centered/bordered
#     header style from BG1, then set the text) ---
$ws.Range("BG1").Copy()
$ws.Range("BH1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("BH1").Value = "Agosto.2021"

# --- Column BH repeats column BG for every existing data row (2-73) ---
for ($r = 2; $r -le 73; $r++) {
    $ws.Range("BH" + $r).Value = $ws.Range("BG" + $r).Value2
}

# --- Row 74 is revised in the new vintage ---
$ws.Range("BH74").Value = 7787

# --- New row 75: new observation period "01-04-2021" (stored as plain text,
#     matching the rest of column A) plus its single data point in BH ---
$ws.Range("Z200").Formula = "=""01-04-2021"""
$ws.Range("Z200").Copy()
$ws.Range("A75").PasteSpecial(-4163)   # xlPasteValues (keeps it text, no date autoconversion)
$ws.Range("Z200").ClearContents()

$ws.Range("BH75").Value = 8218
